# Refactoring to split Supervisor and User requests - trying requests twice if needed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of data (A2:B11) need to be reordered. Use a temporary helper
# column with the desired rank of each row, then sort the whole A:C block
# on that helper column - this moves each row (value + style) as a single
# unit, the same way a real drag/sort in Excel would, instead of copying
# cell-by-cell.
$ranks = @(10, 1, 8, 2, 3, 4, 5, 9, 6, 7)
$r = 2
foreach ($rank in $ranks) {
    $ws.Cells.Item($r, 3).Value = $rank
    $r = $r + 1
}

$sortRange = $ws.Range("A2:C11")
$sortRange.Sort($ws.Range("C2:C11"), 1)

$ws.Range("C2:C11").Clear()

$ws.Range("B15").Select()
